$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), copying the formatting
# from the existing header cell H1 (bold, centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-41: I column and J column values
$data = @(
  @(2,5,5),
  @(3,8,8),
  @(4,6,6),
  @(5,8,8),
  @(6,8,8),
  @(7,7,8),
  @(8,6,6),
  @(9,8,9),
  @(10,8,8),
  @(11,5,5),
  @(12,9,9),
  @(13,10,11),
  @(14,7,8),
  @(15,9,9),
  @(16,5,5),
  @(17,6,7),
  @(18,6,6),
  @(19,6,6),
  @(20,8,8),
  @(21,5,5),
  @(22,9,9),
  @(23,7,8),
  @(24,8,8),
  @(25,4,5),
  @(26,3,4),
  @(27,9,9),
  @(28,7,8),
  @(29,7,7),
  @(30,8,8),
  @(31,9,9),
  @(32,7,7),
  @(33,12,12),
  @(34,8,8),
  @(35,9,9),
  @(36,9,9),
  @(37,4,4),
  @(38,9,9),
  @(39,2,2),
  @(40,5,5),
  @(41,5,5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
